$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values for rows 2, 5, 6, 7
$ws.Range("F2").Value = 0
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -2
